$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.898.48'
$ws.Range('E2').Value = '  -3.34%  '
$ws.Range('D3').Value = '1.827.52'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9974'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '278.10'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -7.32%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5108'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -4.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3476'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -7.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.64'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06796'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.90'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -7.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8098'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -8.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07813'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('D14').Value = '1.813.42'
$ws.Range('E14').Value = '  -3.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.076'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.03'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9963'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.16'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008065'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -4.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9981'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '25.933.55'
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.769'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.02'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.188'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.358'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.99'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.664'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.20'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '109.39'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.334'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -7.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.295'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.10%  '
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04863'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.168'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7300'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -9.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.859'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.182'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.400'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -10.87%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01852'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.72%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5124'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -16.25%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9481'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -10.97%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '116.94'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.208'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.16%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.010'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -8.73%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9973'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1363'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.26%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4490'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -14.82%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.272'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -6.74%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.24'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05917'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.496'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -9.11%  '
